## Add the "Missing Values" worksheet (mirrors the existing "Normalisation"
## sheet layout) and make it the active tab, as described by the commit
## "Done the same with Missing Values / added time taken and code for
## missing values in GOMS".

$wb = $excel.ActiveWorkbook

# --- create the new worksheet after the last existing sheet ----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Missing Values"

# --- table content -----------------------------------------------------
#  row 1: header (Action / Time / Code)
#  rows 2-6: task steps
#  row 7: overall summary
$data = @(
    @("Action",               "Time",   "Code"),
    @("Upload CSV",           "5 min",  "df = pd.read_csv('file.csv')"),
    @("Check Missing Values", "1 min",  "df.isnull().sum()"),
    @("Choose Strategy",      "5 min",  "Decide on filling with mean, median, or mode"),
    @("Fill Missing Values",  "1 min",  "df.fillna(df.mean(), inplace=True)"),
    @("Verify",               "1 min",  "df.isnull().sum()"),
    @("Overall",              "13 min", "")
)

$lastRow = $data.Count
for ($r = 0; $r -lt $data.Count; $r++) {
    $rowNum = $r + 1
    for ($c = 0; $c -lt 3; $c++) {
        $text = $data[$r][$c]
        if ($text -ne "") {
            $ws.Cells.Item($rowNum, $c + 1).Value = $text
        }
    }

    # Row 1 (header) is bold everywhere; row $lastRow ("Overall") is bold
    # only in columns A:B (its C cell is left empty/plain), the rest of the
    # rows use the regular weight - mirrors the "Normalisation" sheet.
    $isBoldRow = ($rowNum -eq 1)
    $boldAB = ($rowNum -eq 1 -or $rowNum -eq $lastRow)

    $cellA = $ws.Range("A" + $rowNum)
    $cellB = $ws.Range("B" + $rowNum)
    $cellC = $ws.Range("C" + $rowNum)

    foreach ($cell in @($cellA, $cellB)) {
        $cell.Font.Size = 13
        $cell.Font.Bold = $boldAB
        $cell.Font.Color = 0
    }

    $cellC.Font.Size = 13
    $cellC.Font.Bold = $isBoldRow
    $cellC.Font.Color = 0

    $ws.Rows.Item($rowNum).RowHeight = 17
}

# --- selection / view state ---------------------------------------------
[void]$ws.Range("A1:C7").Select()
[void]$ws.Activate()

Write-Output "Missing Values sheet added with $($data.Count) rows"
